# Update the build timestamp embedded in the version strings throughout the
# workbook, replacing the old build time with the new one.

$oldText = "February 03 2026 17.29.55 EST"
$newText = "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value2
        if ($v -ne $null -and $v.GetType().Name -eq "String" -and $v.Contains($oldText)) {
            $cell.Value2 = $v.Replace($oldText, $newText)
        }
    }
}
